# Insert a new row at position 176; this shifts the existing rows 176-214
# down to 177-215, and updates the sheet dimension automatically.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows.Item(176).Insert()

# Populate the new row 176 with the weekly market-report entry that was
# added ahead of the previously-most-recent record (now row 177).
$ws.Cells.Item(176, 1).Value = 10
$ws.Cells.Item(176, 2).Value = "Vega Modelo de Temuco"
$ws.Cells.Item(176, 3).Value = "La Araucanía"
$ws.Cells.Item(176, 4).Value = 44889
$ws.Cells.Item(176, 5).Value = 9
$ws.Cells.Item(176, 6).Value = 100114007
$ws.Cells.Item(176, 7).Value = "Jengibre"
$ws.Cells.Item(176, 8).Value = "Sin especificar"
$ws.Cells.Item(176, 9).Value = "Primera"
$ws.Cells.Item(176, 10).Value = 6
$ws.Cells.Item(176, 11).Value = 20000
$ws.Cells.Item(176, 12).Value = 20000
$ws.Cells.Item(176, 13).Value = 20000
$ws.Cells.Item(176, 14).Value = "$/caja 13 kilos"
$ws.Cells.Item(176, 15).Value = "Perú"
$ws.Cells.Item(176, 16).Value = 1538
$ws.Cells.Item(176, 17).Value = 13
$ws.Cells.Item(176, 18).Value = "Hortaliza"
